$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "Qtd_Nós"
$ws.Range("B1").Value = "Ativos"
$ws.Range("C1").Value = "Distancia"
$ws.Range("D1").Value = "Tempo"

# Update data row 2 with the new consolidated values
$ws.Range("A2").Value = 42
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = 3144
$ws.Range("D2").Value = 0.02959275245666504

# Remove leftover column E (old "Tempo" column no longer used)
$ws.Range("E1:E11").ClearContents()

# Remove old rows 3 through 11, which are no longer part of the data
$ws.Range("A3:E11").ClearContents()

$wb.Save()
